# Append 4 new match rows (196-199) to Sheet1, continuing the existing
# fixtures table (League Stats base XGB minute export - Round 20 upload).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (195) down onto the
# four new rows so column A keeps the bold/centered/bordered style used
# for the row-index column throughout the table.
$ws.Range("A195:O195").Copy()
$ws.Range("A196:O199").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New fixture rows, in column order:
# A: index, B: homeTeam, C: awayTeam, D: homeScore, E: awayScore,
# F: homeXg, G: awayXg, H: homeXgPred, I: awayXgPred, J: homePred,
# K: awayPred, L: homeDiff, M: awayDiff, N: totalDiff, O: goalDiff
$data = @(
    @(194, "Monza",    "Fiorentina", 2, 1, 0.67, 1.36, 0.95, 1.75, 0, 1, 0.28, 0.39, 0.67, 2),
    @(195, "Atalanta",  "Juventus",  1, 1, 1.43, 1.56, 2.16, 1.59, 0, 0, 0.73, 0.03, 0.76, 2),
    @(196, "Como",      "Milan",     1, 2, 1.43, 1.35, 1.44, 1.25, 0, 0, 0.01, 0.10, 0.11, 3),
    @(197, "Inter",     "Bologna",   2, 2, 1.79, 0.64, 1.57, 0.87, 1, 0, 0.22, 0.23, 0.45, 3)
)

$r = 196
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}
